$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33 (pushes existing rows 33-40 down to 34-41) so the
# new entry can be slotted in alphabetically with the rest of the list
$ws.Rows("33:33").Insert()

# Populate the new row with the foreclosure-response entry
$ws.Range("A33").Value = "Respond to a mortgage foreclosure complaint"
$ws.Range("B33").Value = "https://www.illinoislegalaid.org/legal-information/respond-mortgage-foreclosure-complaint"
$ws.Range("B33").Style = "Hyperlink"

# The row insert does not automatically re-target the existing hyperlink
# objects onto their shifted cells, so rebuild the hyperlink list to match
# the rows each entry now lives on.
$ws.Hyperlinks.Delete()

$links = @(
    @{ Row = 2; Url = "https://www.illinoislegalaid.org/legal-information/appearance" }
    @{ Row = 15; Url = "https://www.illinoislegalaid.org/legal-information/fee-waiver" }
    @{ Row = 6; Url = "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter" }
    @{ Row = 4; Url = "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting" }
    @{ Row = 12; Url = "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand" }
    @{ Row = 36; Url = "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter" }
    @{ Row = 17; Url = "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr" }
    @{ Row = 39; Url = "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter" }
    @{ Row = 31; Url = "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter" }
    @{ Row = 8; Url = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court" }
    @{ Row = 9; Url = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court" }
    @{ Row = 10; Url = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court" }
    @{ Row = 32; Url = "https://www.illinoislegalaid.org/legal-information/respond-lawsuit" }
    @{ Row = 41; Url = "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap" }
    @{ Row = 18; Url = "https://www.illinoislegalaid.org/legal-information/interpreter-request" }
    @{ Row = 20; Url = "https://www.illinoislegalaid.org/legal-information/motion" }
    @{ Row = 40; Url = "https://www.illinoislegalaid.org/legal-information/transfer-death-instrument-or-todi" }
    @{ Row = 27; Url = "https://www.illinoislegalaid.org/legal-information/power-attorney-agent-resign-letter" }
    @{ Row = 28; Url = "https://www.illinoislegalaid.org/legal-information/power-attorney-revocation" }
    @{ Row = 26; Url = "https://www.illinoislegalaid.org/legal-information/power-attorney-property" }
    @{ Row = 25; Url = "https://www.illinoislegalaid.org/legal-information/power-attorney-health-care" }
    @{ Row = 23; Url = "https://www.illinoislegalaid.org/legal-information/order-protection" }
    @{ Row = 21; Url = "https://www.illinoislegalaid.org/legal-information/name-change-adult" }
    @{ Row = 5; Url = "https://www.illinoislegalaid.org/legal-information/cannabis-expungement" }
    @{ Row = 11; Url = "https://www.illinoislegalaid.org/legal-information/emergency-order-protection-cook-county" }
    @{ Row = 37; Url = "https://www.illinoislegalaid.org/legal-information/short-term-guardian-appointment" }
    @{ Row = 30; Url = "https://www.illinoislegalaid.org/legal-information/remove-eviction-public-record" }
    @{ Row = 34; Url = "https://www.illinoislegalaid.org/legal-information/respond-eviction" }
    @{ Row = 38; Url = "https://www.illinoislegalaid.org/legal-information/small-claims-complaint" }
    @{ Row = 35; Url = "https://www.illinoislegalaid.org/legal-information/security-deposit-complaint" }
    @{ Row = 7; Url = "https://www.illinoislegalaid.org/legal-information/criminal-court-fee-waiver" }
    @{ Row = 16; Url = "https://www.illinoislegalaid.org/legal-information/financial-affidavit" }
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range("B$($link.Row)"), $link.Url) | Out-Null
}

# Reflect the updated scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B38").Select()
